# KP-7577: generating worksheet name for filtered/cascading questions.
# Rename the "Translations question" sheet so its tab name matches the
# auto-generated worksheet-name pattern used for filtered/cascading
# questions ("@@_question").

$wb = $excel.ActiveWorkbook

$targetSheet = $null
foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Translations question") {
        $targetSheet = $ws
    }
}

if ($targetSheet -eq $null) {
    # Fall back to the known position (second sheet) if the name was
    # already changed or differs for some reason.
    $targetSheet = $wb.Worksheets.Item(2)
}

$targetSheet.Name = "@@_question"

Write-Host "Sheet2 name is now:" $targetSheet.Name
